# Fruta / hortaliza, semanal
# Insert two new weekly price records (Murcott - Primera / Segunda) right before the
# existing row 118 entry, shifting the remaining rows (old 118-129) down to 120-131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 118:119 - everything currently at row 118 and below moves down by 2.
$ws.Range("A118:A119").EntireRow.Insert()

# --- New row 118: Murcott / Primera ---
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = 44461
$ws.Range("D118").NumberFormat = $ws.Range("D120").NumberFormat
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100102
$ws.Range("H118").Value = "Cítricos"
$ws.Range("I118").Value = 100102004
$ws.Range("J118").Value = "Mandarina"
$ws.Range("K118").Value = "Murcott"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 240
$ws.Range("N118").Value = 5500
$ws.Range("O118").Value = 6000
$ws.Range("P118").Value = 5750
$ws.Range("Q118").Value = "$/bandeja 10 kilos"
$ws.Range("R118").Value = "Provincia de Limarí"
$ws.Range("S118").Value = 575
$ws.Range("T118").Value = 10

# --- New row 119: Murcott / Segunda ---
$ws.Range("A119").Value = 7
$ws.Range("B119").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C119").Value = "Ñuble"
$ws.Range("D119").Value = 44461
$ws.Range("D119").NumberFormat = $ws.Range("D120").NumberFormat
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100102
$ws.Range("H119").Value = "Cítricos"
$ws.Range("I119").Value = 100102004
$ws.Range("J119").Value = "Mandarina"
$ws.Range("K119").Value = "Murcott"
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 90
$ws.Range("N119").Value = 5000
$ws.Range("O119").Value = 5000
$ws.Range("P119").Value = 5000
$ws.Range("Q119").Value = "$/bandeja 10 kilos"
$ws.Range("R119").Value = "Provincia de Limarí"
$ws.Range("S119").Value = 500
$ws.Range("T119").Value = 10
